$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2511.5557
$ws.Range("I62").Value = 1068.1666
$ws.Range("J62").Value = 3233.25
$ws.Range("K62").Value = 1068.1666
$ws.Range("L62").Value = 3233.25
$ws.Range("M62").Value = -444.1666
$ws.Range("N62").Value = -4481.25
$ws.Range("H64").Value = 45600
$ws.Range("I64").Value = 4583.3335
$ws.Range("J64").Value = 86616.664
$ws.Range("K64").Value = 4583.3335
$ws.Range("L64").Value = 86616.664
$ws.Range("M64").Value = -4335.3335
$ws.Range("N64").Value = -87112.664
$ws.Range("H65").Value = 2511.5557
$ws.Range("I65").Value = 1068.1666
$ws.Range("J65").Value = 3233.25
$ws.Range("K65").Value = 5340.833000000001
$ws.Range("L65").Value = 16166.25
$ws.Range("M65").Value = -2220.833000000001
$ws.Range("N65").Value = -22406.25
$ws.Range("H67").Value = 45600
$ws.Range("I67").Value = 4583.3335
$ws.Range("J67").Value = 86616.664
$ws.Range("K67").Value = 4583.3335
$ws.Range("L67").Value = 86616.664
$ws.Range("M67").Value = -3725.3335
$ws.Range("N67").Value = -88332.664
$ws.Range("H74").Value = 4760
$ws.Range("I74").Value = 4933.3335
$ws.Range("J74").Value = 4500
$ws.Range("K74").Value = 4933.3335
$ws.Range("L74").Value = 4500
$ws.Range("M74").Value = -3997.3335
$ws.Range("N74").Value = -6372
$ws.Range("H76").Value = 3280
$ws.Range("I76").Value = 3466.6667
$ws.Range("J76").Value = 3000
$ws.Range("K76").Value = 3466.6667
$ws.Range("L76").Value = 3000
$ws.Range("M76").Value = -3151.6667
$ws.Range("N76").Value = -3630
$ws.Range("H77").Value = 4760
$ws.Range("I77").Value = 4933.3335
$ws.Range("J77").Value = 4500
$ws.Range("K77").Value = 24666.6675
$ws.Range("L77").Value = 22500
$ws.Range("M77").Value = -19986.6675
$ws.Range("N77").Value = -31860
$ws.Range("H79").Value = 3280
$ws.Range("I79").Value = 3466.6667
$ws.Range("J79").Value = 3000
$ws.Range("K79").Value = 3466.6667
$ws.Range("L79").Value = 3000
$ws.Range("M79").Value = -2374.6667
$ws.Range("N79").Value = -5184
$ws.Range("H125").Value = 2123.3845
$ws.Range("I125").Value = 1954
$ws.Range("J125").Value = 2394.4
$ws.Range("K125").Value = 17586
$ws.Range("L125").Value = 21549.6
$ws.Range("M125").Value = -15126
$ws.Range("N125").Value = -26469.6
$ws.Range("H129").Value = 1334.7693
$ws.Range("I129").Value = 318.75
$ws.Range("J129").Value = 1519.5
$ws.Range("K129").Value = 956.25
$ws.Range("L129").Value = 4558.5
$ws.Range("M129").Value = 4043.75
$ws.Range("N129").Value = -14558.5
$ws.Range("H137").Value = 33335094
$ws.Range("I137").Value = 6174344.5
$ws.Range("K137").Value = 18523033.5
$ws.Range("M137").Value = -18520483.5
$ws.Range("H138").Value = 1430795.8
$ws.Range("I138").Value = 2129081.8
$ws.Range("K138").Value = 6387245.399999999
$ws.Range("M138").Value = -6382105.399999999
$ws.Range("H141").Value = 838.25714
$ws.Range("I141").Value = 685.13794
$ws.Range("J141").Value = 1578.3334
$ws.Range("K141").Value = 2055.41382
$ws.Range("L141").Value = 4735.0002
$ws.Range("M141").Value = 3124.58618
$ws.Range("N141").Value = -15095.0002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 875.0833
$ws.Range("I45").Value = 779.1111
$ws.Range("K45").Value = 779.1111
$ws.Range("M45").Value = -402.1111
$ws.Range("H74").Value = 5719039.5
$ws.Range("I74").Value = 10000729
$ws.Range("K74").Value = 10000729
$ws.Range("M74").Value = -9999855
$ws.Range("H77").Value = 5719039.5
$ws.Range("I77").Value = 10000729
$ws.Range("K77").Value = 50003645
$ws.Range("M77").Value = -49999277
$ws.Range("H102").Value = 1681.4667
$ws.Range("I102").Value = 1623
$ws.Range("J102").Value = 2500
$ws.Range("K102").Value = 1623
$ws.Range("L102").Value = 2500
$ws.Range("M102").Value = -1
$ws.Range("N102").Value = -5744
$ws.Range("H110").Value = 1395
$ws.Range("I110").Value = 1280.25
$ws.Range("J110").Value = 1578.6
$ws.Range("K110").Value = 1280.25
$ws.Range("L110").Value = 1578.6
$ws.Range("M110").Value = 764.75
$ws.Range("N110").Value = -5668.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1384.4736
$ws.Range("I20").Value = 1563.3636
$ws.Range("J20").Value = 1138.5
$ws.Range("K20").Value = 1563.3636
$ws.Range("L20").Value = 1138.5
$ws.Range("M20").Value = -1316.3636
$ws.Range("N20").Value = -1632.5
$ws.Range("H105").Value = 2861
$ws.Range("I105").Value = 3842
$ws.Range("K105").Value = 3842
$ws.Range("M105").Value = -2095

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1523.871
$ws.Range("I31").Value = 1428.2106
$ws.Range("J31").Value = 1675.3334
$ws.Range("K31").Value = 1428.2106
$ws.Range("L31").Value = 1675.3334
$ws.Range("M31").Value = -1133.2106
$ws.Range("N31").Value = -2265.3334
$ws.Range("H34").Value = 1523.871
$ws.Range("I34").Value = 1428.2106
$ws.Range("J34").Value = 1675.3334
$ws.Range("K34").Value = 1428.2106
$ws.Range("L34").Value = 1675.3334
$ws.Range("M34").Value = -1226.2106
$ws.Range("N34").Value = -2079.3334

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 13255
$ws.Range("I80").Value = 4268.3335
$ws.Range("J80").Value = 16625
$ws.Range("K80").Value = 4268.3335
$ws.Range("L80").Value = 16625
$ws.Range("M80").Value = -3270.3335
$ws.Range("N80").Value = -18621
$ws.Range("H83").Value = 13255
$ws.Range("I83").Value = 4268.3335
$ws.Range("J83").Value = 16625
$ws.Range("K83").Value = 21341.6675
$ws.Range("L83").Value = 83125
$ws.Range("M83").Value = -16349.6675
$ws.Range("N83").Value = -93109

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2707.0833
$ws.Range("I40").Value = 2650.4
$ws.Range("J40").Value = 2990.5
$ws.Range("K40").Value = 2650.4
$ws.Range("L40").Value = 2990.5
$ws.Range("M40").Value = -2514.4
$ws.Range("N40").Value = -3262.5
$ws.Range("H61").Value = 1356.579
$ws.Range("I61").Value = 1209.7273
$ws.Range("J61").Value = 1558.5
$ws.Range("K61").Value = 1209.7273
$ws.Range("L61").Value = 1558.5
$ws.Range("M61").Value = -1007.7273
$ws.Range("N61").Value = -1962.5
$ws.Range("H105").Value = 41407.5
$ws.Range("J105").Value = 41407.5
$ws.Range("L105").Value = 41407.5
$ws.Range("N105").Value = -48395.5
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()
$ws.Range("H113").Value = 1356.579
$ws.Range("I113").Value = 1209.7273
$ws.Range("J113").Value = 1558.5
$ws.Range("K113").Value = 1209.7273
$ws.Range("L113").Value = 1558.5
$ws.Range("M113").Value = 960.2727
$ws.Range("N113").Value = -5898.5

